$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("G6").Value = 2.35
$ws.Range("I6").Value = 3.6
$ws.Range("W6").Value = 5.5
$ws.Range("X6").Value = 9.5
$ws.Range("Z6").Value = 23
$ws.Range("AE6").Value = 21
$ws.Range("AI6").Value = 15
$ws.Range("AX6").Value = 5

# Row 7
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.7

# Row 9
$ws.Range("G9").Value = 2.63
$ws.Range("I9").Value = 3
$ws.Range("M9").Value = 1.17
$ws.Range("N9").Value = 5
$ws.Range("Z9").Value = 29
$ws.Range("AE9").Value = 21
$ws.Range("AO9").Value = 19
$ws.Range("BA9").Value = 67

# Row 17
$ws.Range("Q17").Value = 1.95
$ws.Range("R17").Value = 1.8
$ws.Range("S17").Value = 1.4

# Row 46
$ws.Range("G46").Value = 2.25
$ws.Range("H46").Value = 3.45
$ws.Range("J46").Value = 2.85
$ws.Range("K46").Value = 2.12
$ws.Range("L46").Value = 3.35
$ws.Range("O46").Value = 1.2
$ws.Range("P46").Value = 3.6
$ws.Range("W46").Value = 10
$ws.Range("X46").Value = 12.5
$ws.Range("Z46").Value = 24
$ws.Range("AA46").Value = 16.5
$ws.Range("AB46").Value = 22
$ws.Range("AC46").Value = 13
$ws.Range("AD46").Value = 6.9
$ws.Range("AE46").Value = 11.75
$ws.Range("AH46").Value = 11.5
$ws.Range("AI46").Value = 16.5
$ws.Range("AJ46").Value = 10.25
$ws.Range("AM46").Value = 25
$ws.Range("AN46").Value = 4.3
$ws.Range("AO46").Value = 11.75
$ws.Range("AP46").Value = 18.5
$ws.Range("AQ46").Value = 45
$ws.Range("AR46").Value = 75
$ws.Range("AT46").Value = 2.67
$ws.Range("AU46").Value = 6.6
$ws.Range("AX46").Value = 4.85
$ws.Range("AY46").Value = 15
$ws.Range("AZ46").Value = 20
$ws.Range("BB46").Value = 90
$ws.Range("BC46").Value = 250

# Row 47
$ws.Range("N47").Value = 13.3

# Row 61
$ws.Range("G61").Value = 1.75
$ws.Range("H61").Value = 3.3
$ws.Range("I61").Value = 4.2
$ws.Range("J61").Value = 2.5
$ws.Range("K61").Value = 2.05
$ws.Range("O61").Value = 1.4
$ws.Range("P61").Value = 2.75
$ws.Range("Q61").Value = 2.3
$ws.Range("R61").Value = 1.6
$ws.Range("AD61").Value = 7
$ws.Range("AH61").Value = 10
$ws.Range("AP61").Value = 23
$ws.Range("AQ61").Value = 34

# Row 62
$ws.Range("G62").Value = 2.35
$ws.Range("I62").Value = 2.9
$ws.Range("J62").Value = 3.25
$ws.Range("L62").Value = 3.75
$ws.Range("U62").Value = 2.1
$ws.Range("V62").Value = 1.67
$ws.Range("W62").Value = 6.5
$ws.Range("X62").Value = 10
$ws.Range("Y62").Value = 10
$ws.Range("Z62").Value = 23
$ws.Range("AJ62").Value = 12
$ws.Range("AK62").Value = 34
$ws.Range("AL62").Value = 29
$ws.Range("AN62").Value = 4.33
$ws.Range("AX62").Value = 5
$ws.Range("AY62").Value = 19
$ws.Range("AZ62").Value = 34
$ws.Range("BA62").Value = 67
$ws.Range("BC62").Value = 301

# Row 63
$ws.Range("G63").Value = 2.35
$ws.Range("J63").Value = 3.2
$ws.Range("M63").Value = 1.07
$ws.Range("N63").Value = 9
$ws.Range("O63").Value = 1.36
$ws.Range("P63").Value = 3
$ws.Range("Q63").Value = 2.1
$ws.Range("R63").Value = 1.7
$ws.Range("Y63").Value = 10
$ws.Range("AB63").Value = 34
$ws.Range("AC63").Value = 9
$ws.Range("AE63").Value = 17
$ws.Range("AH63").Value = 8
$ws.Range("AN63").Value = 4.5
$ws.Range("AO63").Value = 15
$ws.Range("AP63").Value = 26
$ws.Range("AQ63").Value = 51
$ws.Range("AZ63").Value = 29
